$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row labels: "_old" -> "_FV2210" and "_new" -> "_FV2304" ---
$newHeaders1 = @("Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210","Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210")
$newHeaders2 = @("Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304","Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304")

# Columns A-J (1-10)
for ($i = 0; $i -lt $newHeaders1.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $newHeaders1[$i]
}

# Column K (11) holds "diff" and is left untouched.

# Columns L-U (12-21)
for ($i = 0; $i -lt $newHeaders2.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $newHeaders2[$i]
}

# --- 2. Turn the data range into an Excel Table (ListObject) ---
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U66"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
